$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new manager rows (4, 5, 6) to the list.
# Cells are written in this particular order so that new shared-string
# entries are appended to the shared string table in the same order
# as the target workbook.
$ws.Range("A4").Value = "James"
$ws.Range("A6").Value = "Kelly"
$ws.Range("B4").Value = "T1212121C"
$ws.Range("B5").Value = "S2323232H"
$ws.Range("B6").Value = "T9912834K"
$ws.Range("A5").Value = "Frank"

$ws.Range("C4").Value = 32
$ws.Range("D4").Value = "Married"
$ws.Range("E4").Value = "password"

$ws.Range("C5").Value = 30
$ws.Range("D5").Value = "Single"
$ws.Range("E5").Value = "password"

$ws.Range("C6").Value = 44
$ws.Range("D6").Value = "Married"
$ws.Range("E6").Value = "password"

# Match the final selection in the saved workbook
$ws.Range("G6").Select()
